$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.479.82"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.873.72"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "'315.99"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.5081"
$ws.Range("E7").Value = "  -1.06%  "
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "'0.08367"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").Value = "'1.101"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").Value = "'41.75"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "'6.213"
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").Value = "1.872.10"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "'20.39"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "'7.232"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "'1.010"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "'0.00001101"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "'91.20"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'0.06727"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'17.69"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "'1.008"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "'5.925"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").Value = "28.509.44"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "'11.07"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").Value = "'2.235"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").Value = "2.085.14"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("D27").Value = "'161.78"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").Value = "'20.59"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").Value = "'2.370"
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("D30").Value = "'125.65"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").Value = "'0.1043"
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D32").Value = "'1.034"
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("D33").Value = "'5.759"
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").Value = "'3.619"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "'0.02458"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").Value = "'0.06543"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").Value = "'0.2159"
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").Value = "'8.837"
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("D39").Value = "'5.060"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "'1.186"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").Value = "'0.6392"
$ws.Range("E42").Value = "  -1.08%  "
$ws.Range("D43").Value = "'11.09"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "'1.008"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'0.6009"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").Value = "'13.03"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("D47").Value = "'3.692"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "'2.002"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").Value = "'1.214"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").Value = "'121.79"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("E51").Value = "  -11.90%  "
